$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3 ("GDP vs GDI.xlsx") to make room,
# shifting "GDP vs GDI.xlsx" down to row 4.
$ws.Rows.Item(3).Insert()

# Row 2 currently holds the original Bitcoin template title; update it to the
# new "(with extrapolation)" variant.
$ws.Range("A2").Value = "Bitcoin price change is fueled by global monetary growth (with extrapolation).xlsx"

# Row 3 (newly inserted) gets the original Bitcoin template title.
$ws.Range("A3").Value = "Bitcoin price change is fueled by global monetary growth.xlsx"

# Row 4 still holds "GDP vs GDI.xlsx" (shifted down automatically by the insert).

# Add a new row 5 with the additional template title.
$ws.Range("A5").Value = "Other deposit liabilities (ODL) shows where US M2 is heading.xlsx"
